$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Supplier_On_04/01/19-18:48"
$ws.Range("D2").Value = "Req_Supplier_On_04/01/19-18:48"
$ws.Range("E2").Value = "WorkGroup_On_04/01/19-18:48"
